# Correccion al codigo (Cesantias, incapacidades vacias, fecha liquidacion)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ACUMULATIVOS")

# Row 6
$ws.Range("D6").Value = 2.5
$ws.Range("E6").Value = 150
$ws.Range("G6").Value = 330
$ws.Range("I6").Value = 330
$ws.Range("K6").Value = 45308.60041847919

# Row 7
$ws.Range("D7").Value = 2.5
$ws.Range("E7").Value = 150
$ws.Range("G7").Value = 330
$ws.Range("I7").Value = 330
$ws.Range("K7").Value = 45308.6005060795

# Row 9
$ws.Range("D9").Value = 2.5
$ws.Range("E9").Value = 150
$ws.Range("G9").Value = 330
$ws.Range("I9").Value = 330
$ws.Range("K9").Value = 45308.60057855898

# Row 11
$ws.Range("D11").Value = 2.5
$ws.Range("E11").Value = 150
$ws.Range("G11").Value = 330
$ws.Range("I11").Value = 330
$ws.Range("K11").Value = 45308.60065089612

# Row 12
$ws.Range("D12").Value = 2.5
$ws.Range("E12").Value = 150
$ws.Range("G12").Value = 270
$ws.Range("I12").Value = 270
$ws.Range("K12").Value = 45308.60072208913
